$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the original sheet (with its data, formatting and drawing)
# so we end up with a sheet that will hold the imported data ("to_import")
# while keeping the original sheet (soon renamed "Sheet2") around.
$ws1.Copy($null, $ws1)
$toImport = $wb.Worksheets.Item(2)
$toImport.Name = "to_import"

# The original sheet becomes "Sheet2" and is emptied out.
$ws1.Name = "Sheet2"
$ws1.Cells.Clear()
